# Update the "build" timestamp embedded in the version string wherever it
# appears in the workbook (About sheet header/citation text, and the
# build_version column on the data sheet).
#
# Old: "built on January 30 2026 16.19.47 EST"
# New: "built on February 02 2026 12.49.33 EST"

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$wb = $excel.ActiveWorkbook

# --- "About" sheet: A2 (version banner) and A6 (recommended citation) ---
$wsAbout = $wb.Worksheets.Item("About")

$cellA2 = $wsAbout.Range("A2")
$cellA2.Value = ($cellA2.Text -replace [regex]::Escape($oldStamp), $newStamp)

$cellA6 = $wsAbout.Range("A6")
$cellA6.Value = ($cellA6.Text -replace [regex]::Escape($oldStamp), $newStamp)

# --- "Boundaries and methane sources" sheet: build_version column (S2:S21) ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($r = 2; $r -le 21; $r++) {
    $cell = $wsData.Range("S$r")
    $cell.Value = ($cell.Text -replace [regex]::Escape($oldStamp), $newStamp)
}
